$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update row 27 ("سود هر سهم بر اساس آخرین سرمایه") values per new read_price algorithm
$ws.Range("D27").Value = 47
$ws.Range("E27").Value = 82
$ws.Range("F27").Value = 199
$ws.Range("G27").Value = 434
$ws.Range("H27").Value = 834
